{"js": "// Applies four textual edits to the 1828 music/history timeline document:\n//  1. Zoo paragraph: add \" in Regent's Park, London\" before the final period.\n//  2. Paganini/Bianchi paragraph: \"offers\" -> \"agrees\", and the \"if\" clause\n//     becomes a new sentence (\"In return, ...\").\n//  3. Magistrate's court paragraph: \"A\" -> \"Pursuant to the agreement of 28\n//     July, a\", and \"in Vienna \" is inserted before \"awards\".\n//  4. Trailing date line: \"8 April 2016\" -> \"8 June 2016\".\n\nconst body = context.document.body;\n\n// --- 1. Zoo / Regent's Park -------------------------------------------\n{\n  const found = body.search(\n    \"Two years after the founding of the Zoological Society of London, the London Zoo is opened to the public.\",\n    { matchCase: true }\n  );\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\n      \"Two years after the founding of the Zoological Society of London, the London Zoo is opened to the public in Regent\\u2019s Park, London.\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// --- 2. Paganini agrees / In return ------------------------------------\n{\n  const found = body.search(\n    \"Nicol\\u00f2 Paganini (45) offers to give his mistress, Antonia Bianchi, 2,000 scudi if she will leave him and give him custody of their three-year-old son.\",\n    { matchCase: true }\n  );\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\n      \"Nicol\\u00f2 Paganini (45) agrees to give his mistress, Antonia Bianchi, 2,000 scudi.  In return, she will leave him and give him custody of their three-year-old son.\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// --- 3. Magistrate's court / Vienna ------------------------------------\n{\n  const found = body.search(\n    \"A magistrate\\u2019s court awards custody of Achilles Paganini to his father, Nicol\\u00f2 (45).\",\n    { matchCase: true }\n  );\n  found.load(\"text\");\n  await context.sync();\n  if (found.items.length > 0) {\n    found.items[0].insertText(\n      \"Pursuant to the agreement of 28 July, a magistrate\\u2019s court in Vienna awards custody of Achilles Paganini to his father, Nicol\\u00f2 (45).\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// --- 4. 8 April 2016 -> 8 June 2016 ------------------------------------\n// The date line is split across three runs: \"8 \", \"A\", \"pril 2016\". Target\n// the \"A\" and \"pril 2016\" runs individually so the \"8 \" run (and the\n// paragraph's other formatting) is left completely untouched.\n{\n  const line = body.search(\"8 April 2016\", { matchCase: true });\n  await context.sync();\n  if (line.items.length > 0) {\n    const aPart = line.items[0].search(\"A\", { matchCase: true });\n    await context.sync();\n    if (aPart.items.length > 0) {\n      aPart.items[0].insertText(\"June\", Word.InsertLocation.replace);\n      await context.sync();\n    }\n  }\n\n  // Re-search fresh (the previous range is now stale after the edit above).\n  const prilPart = body.search(\"pril 2016\", { matchCase: true });\n  prilPart.load(\"text\");\n  await context.sync();\n  if (prilPart.items.length > 0) {\n    prilPart.items[0].insertText(\" 2016\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n", "ps1": "# Applies four textual edits to the 1828 music/history timeline document:\n#  1. Zoo paragraph: add \" in Regent's Park, London\" before the final period.\n#  2. Paganini/Bianchi paragraph: \"offers\" -> \"agrees\", and the \"if\" clause\n#     becomes a new sentence (\"In return, ...\").\n#  3. Magistrate's court paragraph: \"A\" -> \"Pursuant to the agreement of 28\n#     July, a\", and \"in Vienna \" is inserted before \"awards\".\n#  4. Trailing date line: \"8 April 2016\" -> \"8 June 2016\".\n\n$d = $word.ActiveDocument\n\n# --- 1. Zoo / Regent's Park ---------------------------------------------\n$rngZoo = $d.Content\nif ($rngZoo.Find.Execute(\"Two years after the founding of the Zoological Society of London, the London Zoo is opened to the public.\")) {\n    $rngZoo.Text = \"Two years after the founding of the Zoological Society of London, the London Zoo is opened to the public in Regent\u2019s Park, London.\"\n}\n\n# --- 2. Paganini agrees / In return -------------------------------------\n$rngPaganini = $d.Content\nif ($rngPaganini.Find.Execute(\"Nicol\u00f2 Paganini (45) offers to give his mistress, Antonia Bianchi, 2,000 scudi if she will leave him and give him custody of their three-year-old son.\")) {\n    $rngPaganini.Text = \"Nicol\u00f2 Paganini (45) agrees to give his mistress, Antonia Bianchi, 2,000 scudi.  In return, she will leave him and give him custody of their three-year-old son.\"\n}\n\n# --- 3. Magistrate's court / Vienna --------------------------------------\n$rngCourt = $d.Content\nif ($rngCourt.Find.Execute(\"A magistrate\u2019s court awards custody of Achilles Paganini to his father, Nicol\u00f2 (45).\")) {\n    $rngCourt.Text = \"Pursuant to the agreement of 28 July, a magistrate\u2019s court in Vienna awards custody of Achilles Paganini to his father, Nicol\u00f2 (45).\"\n}\n\n# --- 4. 8 April 2016 -> 8 June 2016 --------------------------------------\n$rngDate = $d.Content\nif ($rngDate.Find.Execute(\"8 April 2016\")) {\n    $rngDate.Text = \"8 June 2016\"\n}\n"}
